# Updates Price (D) and Volume(1h) (E) columns for rows 2-51 per source diff.
# Columns D/E hold text values (e.g. "26.297.26", "  +0.53%  "), some of which
# look numeric; force the range to Text format before assignment so Excel
# does not auto-convert them to floating point numbers, then restore the
# default "Normal" style so no stray number-format style is left on the cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceVolumeRange = $ws.Range("D2:E51")
$priceVolumeRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.297.26"
$ws.Range("D3").Value = "1.663.78"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  +0.88%  "
$ws.Range("D5").Value = "219.04"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").Value = "0.5340"
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("E8").Value = "  +1.76%  "
$ws.Range("D9").Value = "0.06398"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").Value = "20.58"
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("D11").Value = "0.07829"
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").Value = "4.568"
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("D13").Value = "1.668.57"
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").Value = "1.892.17"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("D15").Value = "0.5522"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").Value = "0.0₅8218"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").Value = "4.683"
$ws.Range("E19").Value = "  +2.14%  "
$ws.Range("D20").Value = "194.44"
$ws.Range("E20").Value = "  +1.95%  "
$ws.Range("D21").Value = "10.21"
$ws.Range("E21").Value = "  +1.51%  "
$ws.Range("D22").Value = "6.038"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("D24").Value = "146.15"
$ws.Range("E24").Value = "  +3.20%  "
$ws.Range("D25").Value = "0.1234"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").Value = "7.184"
$ws.Range("E26").Value = "  -0.75%  "
$ws.Range("D27").Value = "16.10"
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("E28").Value = "  +3.75%  "
$ws.Range("D29").Value = "0.05850"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").Value = "1.283"
$ws.Range("E30").Value = "  +0.86%  "
$ws.Range("D31").Value = "3.616"
$ws.Range("E31").Value = "  +2.63%  "
$ws.Range("D32").Value = "3.281"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("E33").Value = "  +1.54%  "
$ws.Range("D34").Value = "0.9633"
$ws.Range("D35").Value = "2.826"
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("D39").Value = "0.8684"
$ws.Range("E39").Value = "  +2.42%  "
$ws.Range("D40").Value = "5.876"
$ws.Range("E40").Value = "  +1.47%  "
$ws.Range("D41").Value = "1.052.30"
$ws.Range("E41").Value = "  +2.46%  "
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("D43").Value = "104.76"
$ws.Range("E43").Value = "  +2.01%  "
$ws.Range("D44").Value = "1.802.71"
$ws.Range("D45").Value = "57.81"
$ws.Range("E45").Value = "  +1.29%  "
$ws.Range("E46").Value = "  -5.77%  "
$ws.Range("D47").Value = "1.014"
$ws.Range("E47").Value = "  +1.49%  "
$ws.Range("D48").Value = "0.4385"
$ws.Range("E48").Value = "  +2.02%  "
$ws.Range("D49").Value = "8.005"
$ws.Range("E49").Value = "  +2.14%  "
$ws.Range("D50").Value = "0.05164"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").Value = "1.417"
$ws.Range("E51").Value = "  -3.89%  "

$priceVolumeRange.Style = "Normal"
